# Adds ModelSeed-reconstruction columns (G: modelseed_id, H: flags_for_modeling)
# to the media composition sheet, per "updated media composition and
# generated all draft reconstructions".
#
# Cell values are written in the same row-by-row order the author typed
# them in (this keeps the shared-string table's insertion order faithful
# to the original edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "modelseed_id"

$ws.Range("G2").Value  = "cpd00027"
$ws.Range("G3").Value  = "cpd00082"
$ws.Range("G4").Value  = "cpd00158"
$ws.Range("G5").Value  = "cpd00179"
$ws.Range("G6").Value  = "cpd00208"
$ws.Range("G8").Value  = "cpd00003"
$ws.Range("G7").Value  = "cpd08053"
$ws.Range("G9").Value  = "cpd04145"
$ws.Range("G10").Value = "cpd00984"
$ws.Range("G11").Value = "cpd00205,cpd00099"
$ws.Range("G12").Value = "cpd23056"

$ws.Range("H1").Value  = "flags_for_modeling"
$ws.Range("H12").Value = "May have few reactions"

$ws.Range("G13").Value = "cpd00119"
$ws.Range("G14").Value = "cpd00322"
$ws.Range("G15").Value = "cpd00107"
$ws.Range("G16").Value = "cpd00060"
$ws.Range("G17").Value = "cpd00156"
$ws.Range("G18").Value = "cpd00051"
$ws.Range("G19").Value = "cpd00084"
$ws.Range("G20").Value = "cpd00023"
$ws.Range("G21").Value = "cpd00066"
$ws.Range("G22").Value = "cpd00129"
$ws.Range("G23").Value = "cpd00132"
$ws.Range("G24").Value = "cpd00041"
$ws.Range("G25").Value = "cpd00053"
$ws.Range("G26").Value = "cpd00054"
$ws.Range("G27").Value = "cpd00161"
$ws.Range("G28").Value = "cpd00035"
$ws.Range("G29").Value = "cpd00033"
$ws.Range("G30").Value = "cpd00039"
$ws.Range("G31").Value = "cpd00065"
$ws.Range("G32").Value = "cpd00069"
$ws.Range("G33").Value = "cpd00128"
$ws.Range("G34").Value = "cpd00207"
$ws.Range("G35").Value = "cpd00092"
$ws.Range("G36").Value = "cpd00309"
$ws.Range("G37").Value = "cpd00205,cpd00029"
$ws.Range("G38").Value = "cpd00048,cpd10515"
$ws.Range("G39").Value = "cpd00254,cpd00099"
$ws.Range("G40").Value = "cpd00048,cpd00034"
$ws.Range("G41").Value = "cpd00149,cpd00209"
$ws.Range("G42").Value = "cpd00048,cpd00205,cpd24344"
$ws.Range("G43").Value = "cpd00971,cpd03387"
$ws.Range("G44").Value = "cpd00971,cpd15574"

# --- Cells whose fill was explicitly (re-)applied by the author ---
# (shows up as a distinct, explicitly-applied "no fill" style in the xf table)
$flaggedCells = @("G7", "G9", "G10", "G11", "G12", "H12")
foreach ($addr in $flaggedCells) {
    $ws.Range($addr).Interior.ColorIndex = -4142
}

# --- Sheet view: selection as left by the author ---
$ws.Range("C45").Select()
